$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 3): values 2,3,4,5 in R3:U3 ---
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 5

# --- New formula columns R:U for rows 4-30 ---
# R4 is a single (non-shared) formula cell.
$ws.Range("R4").Formula = "=IF(R`$3=`$M4,1,0)"

# Fill R5:U30 with the relative formula first (creates one shared block).
$ws.Range("R5:U30").Formula = "=IF(R`$3=`$M5,1,0)"

# Then fill S4:U19 on top (creates a second shared block that takes
# precedence over the previous one for rows 4-19 in columns S:U).
$ws.Range("S4:U19").Formula = "=IF(S`$3=`$M4,1,0)"

# --- View state: freeze pane scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("T15").Select()

# --- Conditional formatting: add a color-scale rule over R4:U30 and
#     promote it to the top priority (matches the other 3 color-scale
#     rules already on the sheet). ---
$cf = $ws.Range("R4:U30").FormatConditions.AddColorScale(3)
$cf.SetFirstPriority()
